# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.403.04"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.83"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.17"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("E6").Value = "  -1.33%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07602"
$ws.Range("E8").Value = "  +0.60%  "

$ws.Range("E9").Value = "  -1.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.47"
$ws.Range("E10").Value = "  -1.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07742"
$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.846.90"
$ws.Range("E12").Value = "  -6.96%  "

$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001084"
$ws.Range("E14").Value = "  +9.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6797"
$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.72"
$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.100.64"
$ws.Range("E17").Value = "  -7.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.180"
$ws.Range("E18").Value = "  +0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.422.92"
$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.99"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.444"
$ws.Range("E23").Value = "  -1.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.51"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1398"
$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.363"
$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.63"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.466"
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.302"
$ws.Range("E30").Value = "  +4.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05587"
$ws.Range("E31").Value = "  -2.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.103"
$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7094"
$ws.Range("E36").Value = "  -1.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.586"
$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.233.75"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.766"
$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.437"
$ws.Range("E41").Value = "  +5.65%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9052"
$ws.Range("E42").Value = "  -0.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.93"
$ws.Range("E44").Value = "  +0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.07"
$ws.Range("E45").Value = "  -0.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("E46").Value = "  +3.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.171"
$ws.Range("E47").Value = "  +1.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4021"
$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.963"
$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.677"
$ws.Range("E50").Value = "  -1.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1121"
$ws.Range("E51").Value = "  -0.49%  "
